# Update cryptos list data (price and volume(1h) changes, plus some row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their original text formatting
# so values such as "1.00" or "0.0000151" are not coerced into numbers.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '63.050.67'
$ws.Range('E2').Value = '  -4.82%  '
$ws.Range('D3').Value = '3.108.33'
$ws.Range('E3').Value = '  -5.62%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '557.77'
$ws.Range('E5').Value = '  -5.12%  '
$ws.Range('D6').Value = '160.65'
$ws.Range('E6').Value = '  -10.37%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '0.577'
$ws.Range('E8').Value = '  -9.98%  '
$ws.Range('D9').Value = '3.103.55'
$ws.Range('E9').Value = '  -5.70%  '
$ws.Range('D10').Value = '6.71'
$ws.Range('E10').Value = '  -1.78%  '
$ws.Range('D11').Value = '0.114'
$ws.Range('E11').Value = '  -8.86%  '
$ws.Range('D12').Value = '0.374'
$ws.Range('E12').Value = '  -6.63%  '
$ws.Range('D13').Value = '3.648.82'
$ws.Range('E13').Value = '  -5.66%  '
$ws.Range('D14').Value = '0.128'
$ws.Range('E14').Value = '  -1.66%  '
$ws.Range('D15').Value = '63.135.58'
$ws.Range('E15').Value = '  -4.74%  '
$ws.Range('D16').Value = '24.36'
$ws.Range('E16').Value = '  -8.07%  '
$ws.Range('D17').Value = '3.113.98'
$ws.Range('E17').Value = '  -5.77%  '
$ws.Range('D18').Value = '0.0000151'
$ws.Range('E18').Value = '  -7.03%  '
$ws.Range('D19').Value = '394.22'
$ws.Range('E19').Value = '  -7.37%  '
$ws.Range('D20').Value = '12.34'
$ws.Range('E20').Value = '  -4.90%  '
$ws.Range('D21').Value = '5.13'
$ws.Range('E21').Value = '  -6.32%  '
$ws.Range('D22').Value = '6.98'
$ws.Range('E22').Value = '  -4.04%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('D25').Value = '66.63'
$ws.Range('E25').Value = '  -6.42%  '
$ws.Range('E26').Value = '  -3.59%  '
$ws.Range('D27').Value = '0.473'
$ws.Range('E27').Value = '  -7.17%  '
$ws.Range('D28').Value = '0.0₃0999'
$ws.Range('E28').Value = '  -11.99%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '8.51'
$ws.Range('E29').Value = '  -8.71%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').Value = '1.77'
$ws.Range('E32').Value = '  -7.59%  '
$ws.Range('D33').Value = '20.78'
$ws.Range('E33').Value = '  -6.46%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '4.78'
$ws.Range('E34').Value = '  -7.08%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '6.17'
$ws.Range('E35').Value = '  -5.59%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').Value = '1.09'
$ws.Range('E36').Value = '  -7.53%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '150.76'
$ws.Range('E37').Value = '  -5.14%  '
$ws.Range('D38').Value = '1.30'
$ws.Range('E38').Value = '  -8.45%  '
$ws.Range('D39').Value = '2.684.64'
$ws.Range('E39').Value = '  -5.60%  '
$ws.Range('D40').Value = '1.63'
$ws.Range('E40').Value = '  -8.58%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '4.01'
$ws.Range('E41').Value = '  -6.81%  '
$ws.Range('D42').Value = '23.08'
$ws.Range('E42').Value = '  -11.70%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').Value = '38.17'
$ws.Range('E43').Value = '  -3.64%  '
$ws.Range('D44').Value = '0.692'
$ws.Range('E44').Value = '  -7.22%  '
$ws.Range('D45').Value = '0.0601'
$ws.Range('E45').Value = '  -5.70%  '
$ws.Range('D46').Value = '5.46'
$ws.Range('E46').Value = '  -7.12%  '
$ws.Range('D47').Value = '0.0253'
$ws.Range('E47').Value = '  -5.73%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '280.78'
$ws.Range('E48').Value = '  -10.36%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('D50').Value = '20.52'
$ws.Range('E50').Value = '  -9.54%  '
$ws.Range('D51').Value = '0.0968'
$ws.Range('E51').Value = '  -5.38%  '
